$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete one of the old data rows so the table shrinks from 13 to 12
# workers (row shift preserves the special bottom-border style on the
# last row of the table).
$ws.Rows(17).Delete()

# Update totals
$ws.Range("E11").Value = 666340
$ws.Range("C13").Value = 12

# Rewrite the worker table (rows 16-27) with the new account-statement data
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "9147986"
$ws.Range("D16").Value = "YOLFRY JOSE AHUMADA HERRERA"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 908526

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45521649"
$ws.Range("D17").Value = "ALEXANDRA PATRICIA RUIZ MONTIEL"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 908526

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "45550468"
$ws.Range("D18").Value = "LICETH MARIA UPARELA CORRALES"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 908526

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45423586"
$ws.Range("D19").Value = "AMADA ISABEL JIMENEZ BENAVIDES"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 908526

$ws.Range("B20").Value = "CE"
$ws.Range("C20").Value = "302907"
$ws.Range("D20").Value = "EYAL SHATY"
$ws.Range("E20").Value = "2507"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 950000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45456698"
$ws.Range("D21").Value = "CANDELARIA RODRIGUEZ ACOSTA"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "9148231"
$ws.Range("D22").Value = "EDUIN VALENCIA CANTILLO"
$ws.Range("E22").Value = "2507"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 908526

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "73187549"
$ws.Range("D23").Value = "WILLIAM HEREDIA CABRERA"
$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 908526

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1047383240"
$ws.Range("D24").Value = "PIEDAD DEL CARMEN DIAZ GOMEZ"
$ws.Range("E24").Value = "2507"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 908526

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "8854902"
$ws.Range("D25").Value = "JORGE HUMBERT RAMIREZ CURVELO"
$ws.Range("E25").Value = "2507"
$ws.Range("F25").Value = 56940
$ws.Range("G25").Value = 908526

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1002391217"
$ws.Range("D26").Value = "LIZETH PAOLA CASTILLO VASQUEZ"
$ws.Range("E26").Value = "2507"
$ws.Range("F26").Value = 56940
$ws.Range("G26").Value = 908526

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1214217648"
$ws.Range("D27").Value = "ENRIQUE SEGUNDO GUZMAN GUERRERO"
$ws.Range("E27").Value = "2210"
$ws.Range("F27").Value = 40000
$ws.Range("G27").Value = 1000000
